$d = $word.ActiveDocument

# 1. "Hauptteil der Anleitung und in einzelne Punkte unterteilt."
#    -> "Hauptteil der Anleitung und in mehrere einzelne Punkte unterteilt."
$d.Content.Find.Execute("Hauptteil der Anleitung und in einzelne Punkte unterteilt.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Hauptteil der Anleitung und in mehrere einzelne Punkte unterteilt.", 2)

# 2. Fix typo "Sielen" -> "Spielen"
$d.Content.Find.Execute("Haben teilweise ergänzende Anleitungen für Sielen ohne App", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Haben teilweise ergänzende Anleitungen für Spielen ohne App", 2)

# 3. Append " Tablet empfohlen." after "Anmerken, dass nur eine Person die App benötigt."
$d.Content.Find.Execute("Anmerken, dass nur eine Person die App benötigt.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Anmerken, dass nur eine Person die App benötigt. Tablet empfohlen.", 2)

# 4. Remove the "_GoBack" bookmark
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
